# Daily attendance processing - 2025-12-21 10:28:57
# Applies the attendance-recording update to the Session Analysis Results sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Helper: write a literal text string into a cell without Excel's "smart"
# auto-conversion turning percentage-looking / date-looking text into a
# number or date. We evaluate TRIM() of the text in an unused scratch cell
# (far outside the used range) and then paste-special *values only* into the
# destination - this carries over the already-resolved text verbatim and
# keeps the destination cell's existing style/number-format untouched.
# ---------------------------------------------------------------------------
function Set-LiteralText {
    param($range, [string]$text)

    $scratch = $ws.Range("ZZ1000")
    $scratch.Formula = "=TRIM(""" + $text + """)"
    $scratch.Copy()
    $range.PasteSpecial(-4163)   # xlPasteValues
    $scratch.Clear()
}

# ---------------------------------------------------------------------------
# 1) Top summary block (K/L columns) on the first "Class Statistics" table
# ---------------------------------------------------------------------------
$ws.Range("L6").Value = 147      # Recorded Sessions
$ws.Range("L7").Value = 3        # Missing Sessions
Set-LiteralText $ws.Range("L9") "46.2%"   # Coverage %
Set-LiteralText $ws.Range("L10") "72.9%"  # Average Attendance %

# ---------------------------------------------------------------------------
# 2) Rows where "System" was merged into the "Recorded By" list - reorder the
#    names from "System, dnasr281@gmail.com" to "dnasr281@gmail.com, System"
# ---------------------------------------------------------------------------
$gSwapRows = @(8, 9, 10, 34, 35, 36, 60, 61, 62, 86, 87, 88, 112, 113, 114, `
               138, 139, 140, 164, 167, 191, 194, 218, 221, 245, 248, 272, `
               275, 299, 302)

foreach ($r in $gSwapRows) {
    $ws.Cells.Item($r, 7).Value = "dnasr281@gmail.com, System"
}

# ---------------------------------------------------------------------------
# 3) Newly recorded sessions (row 14/40/66/92/118/144): these were
#    "Not Recorded" (pink highlight, style index 4) and are now "Recorded"
#    (plain highlight, style index 2). Copy the formatting from a row that
#    already carries the "Recorded" style (row 2) onto columns A:I, then fill
#    in the recorder, attendance count and status text.
# ---------------------------------------------------------------------------
$newlyRecorded = @{
    14  = "18/26"
    40  = "22/27"
    66  = "18/26"
    92  = "21/27"
    118 = "27/30"
    144 = "20/23"
}

foreach ($r in $newlyRecorded.Keys) {
    $ws.Range("A2:I2").Copy()
    $ws.Range("A$r`:I$r").PasteSpecial(-4122)   # xlPasteFormats

    $ws.Cells.Item($r, 7).Value = "dnasr281@gmail.com"     # Recorded By
    Set-LiteralText $ws.Cells.Item($r, 8) $newlyRecorded[$r]  # Students
    $ws.Cells.Item($r, 9).Value = "Recorded"                # Status
}

$excel.CutCopyMode = 0

# ---------------------------------------------------------------------------
# 4) Group Statistics table (rows 15-20): Recorded/Missing counts and the
#    derived Coverage % / Avg Attendance % figures.
# ---------------------------------------------------------------------------
$groupStats = @{
    15 = @{ O = 12; P = 1; R = "46.2%"; S = "80.8%" }
    16 = @{ O = 13; P = 0; R = "50.0%"; S = "78.6%" }
    17 = @{ O = 13; P = 0; R = "50.0%"; S = "65.1%" }
    18 = @{ O = 13; P = 0; R = "50.0%"; S = "69.2%" }
    19 = @{ O = 13; P = 0; R = "50.0%"; S = "72.3%" }
    20 = @{ O = 12; P = 1; R = "46.2%"; S = "75.4%" }
}

foreach ($r in $groupStats.Keys) {
    $row = $groupStats[$r]
    $ws.Range("O$r").Value = $row.O
    $ws.Range("P$r").Value = $row.P
    Set-LiteralText $ws.Range("R$r") $row.R
    Set-LiteralText $ws.Range("S$r") $row.S
}
